# The row "BNE | Bruttonationaleinkommen | " (row 154) was removed from the
# abbreviation list. All subsequent rows (155-185) shift up by one row, and
# the used range shrinks from A1:C185 to A1:C184.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Delete()
